$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.056.08"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.363.72"
$ws.Range("E3").Value = "  -4.06%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "502.33"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.60"
$ws.Range("E6").Value = "  -3.36%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  -2.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.366.55"
$ws.Range("E9").Value = "  -3.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0986"
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.83"
$ws.Range("E12").Value = "  +4.00%  "
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.783.40"
$ws.Range("E14").Value = "  -4.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.029.58"
$ws.Range("E15").Value = "  -3.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.45"
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.364.16"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.02"
$ws.Range("E19").Value = "  -3.46%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "307.90"
$ws.Range("E20").Value = "  -2.49%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.01"
$ws.Range("E21").Value = "  -2.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.31"
$ws.Range("E22").Value = "  -2.48%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.95"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.370"
$ws.Range("E26").Value = "  -3.06%  "
$ws.Range("E27").Value = "  -6.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.23"
$ws.Range("E28").Value = "  -5.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.67"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0711"
$ws.Range("E30").Value = "  -3.79%  "
$ws.Range("E31").Value = "  -3.07%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.79"
$ws.Range("E33").Value = "  -5.94%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("E35").Value = "  -5.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.64"
$ws.Range("E36").Value = "  -2.58%  "
$ws.Range("E37").Value = "  -5.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.73"
$ws.Range("E38").Value = "  -4.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.23"
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.39"
$ws.Range("E41").Value = "  -5.83%  "
$ws.Range("E42").Value = "  -1.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "128.96"
$ws.Range("E43").Value = "  -5.77%  "
$ws.Range("E44").Value = "  -5.14%  "
$ws.Range("E45").Value = "  -2.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0903"
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "238.78"
$ws.Range("E47").Value = "  -6.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0481"
$ws.Range("E48").Value = "  -2.91%  "
$ws.Range("E49").Value = "  -4.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.96"
$ws.Range("E50").Value = "  -2.59%  "
$ws.Range("E51").Value = "  -1.28%  "

Write-Host "done"